# Update computed result values on each year's sheet ("ADD results from server")
$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell address -> new value
$updates = @{
    "2025" = @{
        "N2" = 7155.075790473336
        "O2" = 6980.325566461754
    }
    "2030" = @{
        "B2" = 5707.815717280662
        "I2" = 44492.05901988943
        "L2" = 66334.06707325629
        "M2" = 21991.42050229464
        "O2" = 12079.40905079305
    }
    "2035" = @{
        "A2" = 2927.360317916481
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15119.05116437215
        "O2" = 14761.1605844395
    }
    "2040" = @{
        "A2" = 2927.360317916481
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15226.1750723931
        "O2" = 14761.1605844395
    }
    "2045" = @{
        "A2" = 6352.985609279765
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15770.90262384333
        "O2" = 17096.62657078825
    }
    "2050" = @{
        "A2" = 6352.985609279765
        "B2" = 7940.887964949257
        "E2" = 67179.99183625776
        "I2" = 59530.75343380851
        "L2" = 66334.06707325629
        "M2" = 25547.11936466757
        "N2" = 15770.90262384333
        "O2" = 17096.62657078825
    }
}

foreach ($sheetName in $updates.Keys) {
    $nameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($nameStr)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range([string]$addr).Value = $cellUpdates[$addr]
    }
}
